# Apply updated dSF (column F) values to specific rows per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -7
$ws.Range("F5").Value = -9
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -5
$ws.Range("F10").Value = -7
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = -2
$ws.Range("F15").Value = -6
$ws.Range("F18").Value = -2
$ws.Range("F20").Value = -1
